# Add a "Price" column (column E) to each Vendor sheet, matching the
# commit "Updated file with prices".

$wb = $excel.ActiveWorkbook

# ---- Vendor 1 -----------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Vendor 1")
$ws1.Range("E1").Value = "Price"
$ws1.Range("E2").Value = 3.14
$ws1.Range("E3").Value = 2.74
$ws1.Range("E4").Value = 1
$ws1.Range("E4").NumberFormat = "0.00"
$ws1.Range("E5").Value = 68
$ws1.Range("C4").Select() | Out-Null

# ---- Vendor 2 -----------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Vendor 2")
$ws2.Range("E1").Value = "Price "
$ws2.Range("E2").Value = 5.5
$ws2.Range("E2").NumberFormat = "0.00"
$ws2.Range("E3").Value = 4.99
$ws2.Range("E4").Value = 3.61
$ws2.Range("E5").Value = 1
$ws2.Range("E5").NumberFormat = "0.00"
$ws2.Range("E5").Select() | Out-Null

# ---- Vendor 3 -------------------------------------------------------------
# Note: Vendor 3's header row never got a "Price" label in the source edit,
# only the price values were filled in - replicated faithfully here.
$ws3 = $wb.Worksheets.Item("Vendor 3")
$ws3.Range("E2").Value = 1
$ws3.Range("E2").NumberFormat = "0.00"
$ws3.Range("E3").Value = 2.29
$ws3.Range("E4").Value = 2.78
$ws3.Range("E5").Value = 104.99

# Vendor 3 ends up the active/selected tab with F15 selected.
$ws3.Activate() | Out-Null
$ws3.Range("F15").Select() | Out-Null
